# Update odds values for row 4 (match: Correcaminos - Atl. Morelia)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = 3.45
$ws.Range("I4").Value = 2.2
$ws.Range("K4").Value = 2.18
$ws.Range("L4").Value = 2.75
$ws.Range("O4").Value = 1.21
$ws.Range("P4").Value = 3.5
$ws.Range("Q4").Value = 1.65
$ws.Range("R4").Value = 1.98
$ws.Range("S4").Value = 1.33
$ws.Range("T4").Value = 3.14
$ws.Range("U4").Value = 1.55
$ws.Range("V4").Value = 2.15
$ws.Range("W4").Value = 11.25
$ws.Range("X4").Value = 17
$ws.Range("Y4").Value = 10.5
$ws.Range("AA4").Value = 22
$ws.Range("AB4").Value = 26
$ws.Range("AC4").Value = 12.5
$ws.Range("AD4").Value = 6.9
$ws.Range("AE4").Value = 12
$ws.Range("AF4").Value = 45
$ws.Range("AH4").Value = 9.5
$ws.Range("AI4").Value = 12
$ws.Range("AJ4").Value = 8.75
$ws.Range("AK4").Value = 22
$ws.Range("AL4").Value = 16.5
$ws.Range("AM4").Value = 23
$ws.Range("AP4").Value = 20
$ws.Range("AS4").Value = 200
$ws.Range("AT4").Value = 2.92
$ws.Range("AU4").Value = 6.6
$ws.Range("AV4").Value = 50
$ws.Range("AX4").Value = 11.25
$ws.Range("AY4").Value = 17.5
$ws.Range("BA4").Value = 70
